# Update impact values in the recipe templates (test scenarios)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value  = 9.749285
$ws.Range("C2").Value  = 1.543815
$ws.Range("D2").Value  = 2.9939

$ws.Range("B3").Value  = 2.367085
$ws.Range("C3").Value  = 0.371455
$ws.Range("D3").Value  = 0.89425

$ws.Range("B5").Value  = 7.402089999999999
$ws.Range("C5").Value  = 0.43149
$ws.Range("D5").Value  = 1.1513

$ws.Range("B7").Value  = 18.63214
$ws.Range("C7").Value  = 0.935825
$ws.Range("D7").Value  = 1.01636

$ws.Range("B8").Value  = 0.590395
$ws.Range("C8").Value  = 0.37427
$ws.Range("D8").Value  = 0.2132

$ws.Range("B9").Value  = 0.17663
$ws.Range("C9").Value  = 0.09973499999999999
$ws.Range("D9").Value  = 0.14478

$ws.Range("B10").Value = 0.129185
$ws.Range("C10").Value = 0.097335
$ws.Range("D10").Value = 0.11326

$ws.Range("B11").Value = 0.08722000000000001
$ws.Range("C11").Value = 0.13545
$ws.Range("D11").Value = 0.093135

$ws.Range("B12").Value = 0.09084095
$ws.Range("C12").Value = 0.00755775
$ws.Range("D12").Value = 0.03401145000000001

$ws.Range("B13").Value = 0.08669285
$ws.Range("C13").Value = 0.00846015
$ws.Range("D13").Value = 0.01939835

$ws.Range("B14").Value = 788.7049999999999
$ws.Range("C14").Value = 455.645
$ws.Range("D14").Value = 575.3099999999999

$ws.Range("B15").Value = 18284.165
$ws.Range("C15").Value = 10832.175
$ws.Range("D15").Value = 12671.74
